$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 568, shifting the existing weekly records (old rows 568-591)
# down by one (to rows 569-592), and populate the newly inserted row with the
# latest weekly observation.
$ws.Rows("568:568").Insert()

$ws.Cells.Item(568,1).Value  = 6
$ws.Cells.Item(568,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(568,3).Value  = "Metropolitana"
$ws.Cells.Item(568,4).Value  = 44939
$ws.Cells.Item(568,5).Value  = 13
$ws.Cells.Item(568,6).Value  = 100112039
$ws.Cells.Item(568,7).Value  = "Ciboulette"
$ws.Cells.Item(568,8).Value  = "Sin especificar"
$ws.Cells.Item(568,9).Value  = "Primera"
$ws.Cells.Item(568,10).Value = 770
$ws.Cells.Item(568,11).Value = 700
$ws.Cells.Item(568,12).Value = 800
$ws.Cells.Item(568,13).Value = 745
$ws.Cells.Item(568,14).Value = "`$/docena de atados"
$ws.Cells.Item(568,15).Value = "Región Metropolitana"
$ws.Cells.Item(568,16).Value = 248
$ws.Cells.Item(568,17).Value = 3
$ws.Cells.Item(568,18).Value = "Hortaliza"
